$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1

$ws.Range("E3").Select()
